$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value (values are forced to Text
# via a leading apostrophe so Excel doesn't reinterpret numeric-looking
# strings like "1.00" or "0.0580" as numbers, then the style is reset to
# "Normal" so no extra cell formatting / number-format styles get added.
$updates = @(
    @{ Cell = 'D2'; Value = '''61.023.25' }
    @{ Cell = 'E2'; Value = '''  +6.96%  ' }
    @{ Cell = 'D3'; Value = '''2.627.02' }
    @{ Cell = 'E3'; Value = '''  +9.20%  ' }
    @{ Cell = 'D4'; Value = '''1.00' }
    @{ Cell = 'E4'; Value = '''  +0.02%  ' }
    @{ Cell = 'D5'; Value = '''509.83' }
    @{ Cell = 'E5'; Value = '''  +4.35%  ' }
    @{ Cell = 'D6'; Value = '''159.20' }
    @{ Cell = 'E6'; Value = '''  +2.92%  ' }
    @{ Cell = 'D7'; Value = '''0.997' }
    @{ Cell = 'E7'; Value = '''  +0.10%  ' }
    @{ Cell = 'E8'; Value = '''  -1.72%  ' }
    @{ Cell = 'D9'; Value = '''2.665.90' }
    @{ Cell = 'E9'; Value = '''  +9.96%  ' }
    @{ Cell = 'E10'; Value = '''  +3.64%  ' }
    @{ Cell = 'E11'; Value = '''  +5.67%  ' }
    @{ Cell = 'E12'; Value = '''  +3.89%  ' }
    @{ Cell = 'E13'; Value = '''  +1.04%  ' }
    @{ Cell = 'D14'; Value = '''3.088.04' }
    @{ Cell = 'E14'; Value = '''  +9.20%  ' }
    @{ Cell = 'D15'; Value = '''60.904.56' }
    @{ Cell = 'E15'; Value = '''  +6.66%  ' }
    @{ Cell = 'D16'; Value = '''21.79' }
    @{ Cell = 'E16'; Value = '''  +5.84%  ' }
    @{ Cell = 'D17'; Value = '''0.0000142' }
    @{ Cell = 'E17'; Value = '''  +6.40%  ' }
    @{ Cell = 'D18'; Value = '''2.659.47' }
    @{ Cell = 'E19'; Value = '''  +1.93%  ' }
    @{ Cell = 'D20'; Value = '''349.09' }
    @{ Cell = 'E20'; Value = '''  +7.42%  ' }
    @{ Cell = 'D21'; Value = '''10.57' }
    @{ Cell = 'E21'; Value = '''  +6.18%  ' }
    @{ Cell = 'D22'; Value = '''6.21' }
    @{ Cell = 'E22'; Value = '''  +4.67%  ' }
    @{ Cell = 'D23'; Value = '''0.997' }
    @{ Cell = 'E23'; Value = '''  +0.01%  ' }
    @{ Cell = 'D24'; Value = '''60.70' }
    @{ Cell = 'E24'; Value = '''  +4.58%  ' }
    @{ Cell = 'E25'; Value = '''  +5.17%  ' }
    @{ Cell = 'D26'; Value = '''2.745.19' }
    @{ Cell = 'E26'; Value = '''  +9.25%  ' }
    @{ Cell = 'D27'; Value = '''0.167' }
    @{ Cell = 'E27'; Value = '''  +3.79%  ' }
    @{ Cell = 'D28'; Value = '''0.991' }
    @{ Cell = 'E28'; Value = '''  -0.45%  ' }
    @{ Cell = 'D29'; Value = '''0.0₃0874' }
    @{ Cell = 'E29'; Value = '''  +11.56%  ' }
    @{ Cell = 'E30'; Value = '''  +4.79%  ' }
    @{ Cell = 'D31'; Value = '''1.00' }
    @{ Cell = 'E31'; Value = '''  +0.18%  ' }
    @{ Cell = 'D33'; Value = '''157.11' }
    @{ Cell = 'E33'; Value = '''  +4.39%  ' }
    @{ Cell = 'E34'; Value = '''  +3.84%  ' }
    @{ Cell = 'D35'; Value = '''5.85' }
    @{ Cell = 'E35'; Value = '''  +10.94%  ' }
    @{ Cell = 'E36'; Value = '''  +8.36%  ' }
    @{ Cell = 'E37'; Value = '''  +5.73%  ' }
    @{ Cell = 'D38'; Value = '''311.80' }
    @{ Cell = 'E38'; Value = '''  +15.57%  ' }
    @{ Cell = 'B39'; Value = '''Fetch.AI' }
    @{ Cell = 'C39'; Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D39'; Value = '''0.864' }
    @{ Cell = 'E39'; Value = '''  +2.74%  ' }
    @{ Cell = 'B40'; Value = '''Stacks' }
    @{ Cell = 'C40'; Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D40'; Value = '''1.51' }
    @{ Cell = 'E40'; Value = '''  +9.85%  ' }
    @{ Cell = 'E41'; Value = '''  +31.59%  ' }
    @{ Cell = 'D43'; Value = '''35.29' }
    @{ Cell = 'E43'; Value = '''  +3.25%  ' }
    @{ Cell = 'D44'; Value = '''0.638' }
    @{ Cell = 'E44'; Value = '''  +7.09%  ' }
    @{ Cell = 'D45'; Value = '''0.0581' }
    @{ Cell = 'E45'; Value = '''  +9.74%  ' }
    @{ Cell = 'E46'; Value = '''  -1.38%  ' }
    @{ Cell = 'D47'; Value = '''0.994' }
    @{ Cell = 'E47'; Value = '''  +0.01%  ' }
    @{ Cell = 'D48'; Value = '''19.70' }
    @{ Cell = 'E48'; Value = '''  +12.92%  ' }
    @{ Cell = 'D49'; Value = '''4.91' }
    @{ Cell = 'E49'; Value = '''  +8.12%  ' }
    @{ Cell = 'D50'; Value = '''2.073.23' }
    @{ Cell = 'E50'; Value = '''  +10.31%  ' }
    @{ Cell = 'D51'; Value = '''0.0237' }
    @{ Cell = 'E51'; Value = '''  +4.00%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
